# Fix typo in Matt's surname: "Hoffmanna" -> "Hoffmann"
$d = $word.ActiveDocument
$d.Content.Find.Execute("Hoffmanna", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hoffmann", 2)
